# Apply the edit described by the commit:
#   "Specify config file on command line. Rearrange columns to keep 'Reports' last."
#
# Observable effects in the workbook:
#   1. The three "Date of Exam" report strings on the Patients sheet get the
#      word "lateral" parenthesised: "MV e' lateral:" -> "MV e' (lateral):"
#   2. The previously-active sheet ("Labs") is no longer the active tab; the
#      "Patients" sheet becomes active/selected, with its selection moved to B5.
#   3. The Labs sheet's own selection (B2) is left as-is, it simply stops being
#      the active tab.

$wb = $excel.ActiveWorkbook

$patients = $wb.Worksheets.Item("Patients")
$labs = $wb.Worksheets.Item("Labs")

# 1. Fix up the report text bodies (B2:B4) on the Patients sheet.
foreach ($r in 2..4) {
    $cell = $patients.Cells.Item($r, 2)
    $text = $cell.Value2
    $cell.Value2 = $text.Replace("lateral:", "(lateral):")
}

# 2. Make "Patients" the active sheet/tab again, with B5 selected.
$patients.Activate()
$patients.Range("B5").Select()

# Leave the Labs sheet's own selection (B2) untouched; it already is B2.
